# ATA da reuniao do dia 20/04/2020
# - remove the empty "Planilha2" worksheet
# - add the meeting minutes for 20/04/2020 into row 16 of Planilha1
# - adjust row heights (15 & 16) and the saved selection to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Fill in the new meeting row (row 16) ---------------------------------
$ws.Range("A16").Value = "20/042020"

$ws.Range("B16").Value = 0.79166666666666663
$ws.Range("B16").NumberFormat = $ws.Range("B15").NumberFormat

$ws.Range("C16").Value = 0.80069444444444438
$ws.Range("C16").NumberFormat = $ws.Range("C15").NumberFormat

$ws.Range("D16").Value = "Stefany Batista, Graziela, Gabriel Bezerra, Yuri Vedovate, Bruno Santana, Raphael Moitinho"

$ws.Range("F16").Value = "1 - o que fazer essa semana   2 - focar nas atividades de cada dupla.                               3 - quinta feira (23/04/2020) ensaio geral para as apresentações.                                   "

# --- Row heights ------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 143.25
$ws.Rows.Item(16).RowHeight = 143.25

# --- Update the active selection -------------------------------------------
$ws.Range("F16").Select() | Out-Null

# --- Remove the unused second worksheet -------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Planilha2").Delete() | Out-Null
